$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "username"
$ws.Range("L1").Value = "pswd"

for ($r = 2; $r -le 9; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 11).Value = $name
    $ws.Cells.Item($r, 12).Value = "abcd"
}

$ws.Range("F2").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("K9").Select()
